$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 12: MTTH=900, Interval=20, with probability formulas mirroring rows 6-11
$ws.Range("C12").Formula = "=1-(0.5)^(F12/E12)"
$ws.Range("D12").Formula = "=F12*(1-0.5^(1/E12))"
$ws.Range("E12").Value = 900
$ws.Range("F12").Value = 20

# Move the active selection to G10, matching the saved view state
$ws.Range("G10").Select()
